$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 1131.6428
$ws.Cells.Item(2, 9).Value = 230.75
$ws.Cells.Item(2, 10).Value = 2332.8333
$ws.Cells.Item(2, 11).Value = 230.75
$ws.Cells.Item(2, 12).Value = 2332.8333
$ws.Cells.Item(2, 13).Value = -117.75
$ws.Cells.Item(2, 14).Value = -2558.8333

# Row 137
$ws.Cells.Item(137, 8).Value = 3160.246
$ws.Cells.Item(137, 9).Value = 3217.818
$ws.Cells.Item(137, 10).Value = 3148.5186
$ws.Cells.Item(137, 11).Value = 9653.454000000002
$ws.Cells.Item(137, 12).Value = 9445.5558
$ws.Cells.Item(137, 13).Value = -7103.454000000002
$ws.Cells.Item(137, 14).Value = -14545.5558

# Row 138
$ws.Cells.Item(138, 8).Value = 3525.4844
$ws.Cells.Item(138, 10).Value = 3663
$ws.Cells.Item(138, 12).Value = 10989
$ws.Cells.Item(138, 14).Value = -21269

# Row 140
$ws.Cells.Item(140, 8).Value = 69189.25
$ws.Cells.Item(140, 10).Value = 69189.25
$ws.Cells.Item(140, 12).Value = 69189.25
$ws.Cells.Item(140, 14).Value = -79549.25

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 13).ClearContents()

# Row 45
$ws.Cells.Item(45, 8).Value = 2731.5
$ws.Cells.Item(45, 9).Value = 2281.7083
$ws.Cells.Item(45, 11).Value = 2281.7083
$ws.Cells.Item(45, 13).Value = -1904.7083

# Row 74
$ws.Cells.Item(74, 8).Value = 19610386
$ws.Cells.Item(74, 9).Value = 23810734
$ws.Cells.Item(74, 10).Value = 8754.666999999999
$ws.Cells.Item(74, 11).Value = 23810734
$ws.Cells.Item(74, 12).Value = 8754.666999999999
$ws.Cells.Item(74, 13).Value = -23809860
$ws.Cells.Item(74, 14).Value = -10502.667

# Row 77
$ws.Cells.Item(77, 8).Value = 19610386
$ws.Cells.Item(77, 9).Value = 23810734
$ws.Cells.Item(77, 10).Value = 8754.666999999999
$ws.Cells.Item(77, 11).Value = 119053670
$ws.Cells.Item(77, 12).Value = 43773.335
$ws.Cells.Item(77, 13).Value = -119049302
$ws.Cells.Item(77, 14).Value = -52509.335

# Row 122
$ws.Cells.Item(122, 8).Value = 3628.3635
$ws.Cells.Item(122, 9).Value = 3303.0908
$ws.Cells.Item(122, 11).Value = 9909.2724
$ws.Cells.Item(122, 13).Value = -7459.2724

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 1146.2858
$ws.Cells.Item(94, 9).Value = 888.3077
$ws.Cells.Item(94, 11).Value = 888.3077
$ws.Cells.Item(94, 13).Value = -437.3077

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1671
$ws.Cells.Item(16, 9).Value = 755.4167
$ws.Cells.Item(16, 10).Value = 5333.3335
$ws.Cells.Item(16, 11).Value = 755.4167
$ws.Cells.Item(16, 12).Value = 5333.3335
$ws.Cells.Item(16, 13).Value = -468.4167
$ws.Cells.Item(16, 14).Value = -5907.3335

# Row 31
$ws.Cells.Item(31, 8).Value = 44677.08
$ws.Cells.Item(31, 9).Value = 2734.158
$ws.Cells.Item(31, 10).Value = 177496.33
$ws.Cells.Item(31, 11).Value = 2734.158
$ws.Cells.Item(31, 12).Value = 177496.33
$ws.Cells.Item(31, 13).Value = -2439.158
$ws.Cells.Item(31, 14).Value = -178086.33

# Row 34
$ws.Cells.Item(34, 8).Value = 44677.08
$ws.Cells.Item(34, 9).Value = 2734.158
$ws.Cells.Item(34, 10).Value = 177496.33
$ws.Cells.Item(34, 11).Value = 2734.158
$ws.Cells.Item(34, 12).Value = 177496.33
$ws.Cells.Item(34, 13).Value = -2532.158
$ws.Cells.Item(34, 14).Value = -177900.33

# Row 52
$ws.Cells.Item(52, 8).Value = 32779.25
$ws.Cells.Item(52, 9).Value = 28708.5
$ws.Cells.Item(52, 10).Value = 36850
$ws.Cells.Item(52, 11).Value = 28708.5
$ws.Cells.Item(52, 12).Value = 36850
$ws.Cells.Item(52, 13).Value = -28414.5
$ws.Cells.Item(52, 14).Value = -37438

# Row 99
$ws.Cells.Item(99, 8).Value = 2307.6924
$ws.Cells.Item(99, 10).Value = 2500
$ws.Cells.Item(99, 12).Value = 2500
$ws.Cells.Item(99, 14).Value = -5496

# Row 113
$ws.Cells.Item(113, 8).Value = 1671
$ws.Cells.Item(113, 9).Value = 755.4167
$ws.Cells.Item(113, 10).Value = 5333.3335
$ws.Cells.Item(113, 11).Value = 755.4167
$ws.Cells.Item(113, 12).Value = 5333.3335
$ws.Cells.Item(113, 13).Value = 1414.5833
$ws.Cells.Item(113, 14).Value = -9673.333500000001

# Row 126
$ws.Cells.Item(126, 8).Value = 2307.6924
$ws.Cells.Item(126, 10).Value = 2500
$ws.Cells.Item(126, 12).Value = 7500
$ws.Cells.Item(126, 14).Value = -12440

# Row 132
$ws.Cells.Item(132, 8).Value = 2878.6365
$ws.Cells.Item(132, 9).Value = 2068.6177
$ws.Cells.Item(132, 10).Value = 5632.7
$ws.Cells.Item(132, 11).Value = 6205.853099999999
$ws.Cells.Item(132, 12).Value = 16898.1
$ws.Cells.Item(132, 13).Value = -3675.853099999999
$ws.Cells.Item(132, 14).Value = -21958.1

# Row 134
$ws.Cells.Item(134, 8).Value = 2068.5
$ws.Cells.Item(134, 9).Value = 1220.381
$ws.Cells.Item(134, 10).Value = 4047.4443
$ws.Cells.Item(134, 11).Value = 3661.143
$ws.Cells.Item(134, 12).Value = 12142.3329
$ws.Cells.Item(134, 13).Value = -1126.143
$ws.Cells.Item(134, 14).Value = -17212.3329

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 6521050.5
$ws.Cells.Item(131, 9).Value = 25000984
$ws.Cells.Item(131, 10).Value = 4421058
$ws.Cells.Item(131, 11).Value = 75002952
$ws.Cells.Item(131, 12).Value = 13263174
$ws.Cells.Item(131, 13).Value = -74997912
$ws.Cells.Item(131, 14).Value = -13273254

$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Cells.Item(36, 8).Value = 2779.25
$ws.Cells.Item(36, 9).Value = 2779.25
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 2779.25
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(36, 14).Value = -2294.25

# Row 113
$ws.Cells.Item(113, 8).Value = 1850.1052
$ws.Cells.Item(113, 10).Value = 1882.3334
$ws.Cells.Item(113, 12).Value = 1882.3334
$ws.Cells.Item(113, 14).Value = -6222.3334

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 9150
$ws.Cells.Item(40, 9).Value = 8387.875
$ws.Cells.Item(40, 11).Value = 8387.875
$ws.Cells.Item(40, 13).Value = -8251.875

# Row 45
$ws.Cells.Item(45, 8).Value = 5000
$ws.Cells.Item(45, 9).Value = 5000
$ws.Cells.Item(45, 11).Value = 5000
$ws.Cells.Item(45, 13).Value = -4593

# Row 58
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 13).ClearContents()

# Row 122
$ws.Cells.Item(122, 8).Value = 158334.39
$ws.Cells.Item(122, 9).Value = 194204.28
$ws.Cells.Item(122, 10).Value = 7680.8
$ws.Cells.Item(122, 11).Value = 582612.84
$ws.Cells.Item(122, 12).Value = 23042.4
$ws.Cells.Item(122, 13).Value = -580162.84
$ws.Cells.Item(122, 14).Value = -27942.4

$ws = $wb.Worksheets.Item("WVR")
# Row 25
$ws.Cells.Item(25, 8).Value = 30000
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 13).ClearContents()

# Row 107
$ws.Cells.Item(107, 8).Value = 1650.5
$ws.Cells.Item(107, 9).Value = 1509
$ws.Cells.Item(107, 11).Value = 4527
$ws.Cells.Item(107, 13).Value = -2607

# Row 113
$ws.Cells.Item(113, 8).Value = 358.4091
$ws.Cells.Item(113, 9).Value = 319.27777
$ws.Cells.Item(113, 11).Value = 957.83331
$ws.Cells.Item(113, 13).Value = 1212.16669

# Row 122
$ws.Cells.Item(122, 8).Value = 3870.8
$ws.Cells.Item(122, 9).Value = 3870.8
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value = -9162.400000000001
